# Adicionados balancos concatenados em uma unica planilha.
# Adds BF:BH columns (31/12/2023, 31/03/2024, 30/06/2024) to the PRIO3 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy formatting (bold/border/center) from BE1, then set the new period labels.
$ws.Range("BE1").Copy($ws.Range("BF1:BH1"))
$ws.Range("BF1").Value = "31/12/2023"
$ws.Range("BG1").Value = "31/03/2024"
$ws.Range("BH1").Value = "30/06/2024"

# Row 2
$ws.Range("BF2").Value = 28310931.456
$ws.Range("BG2").Value = 30779117.568
$ws.Range("BH2").Value = 37293928.448
# Row 3
$ws.Range("BF3").Value = 5294874.112
$ws.Range("BG3").Value = 6790404.096
$ws.Range("BH3").Value = 10549708.8
# Row 4
$ws.Range("BF4").Value = 2335403.008
$ws.Range("BG4").Value = 4030097.92
$ws.Range("BH4").Value = 6447834.112
# Row 5
$ws.Range("BF5").Value = 0
$ws.Range("BG5").Value = 0
$ws.Range("BH5").Value = 0
# Row 6
$ws.Range("BF6").Value = 1743490.944
$ws.Range("BG6").Value = 1431728
$ws.Range("BH6").Value = 2093881.984
# Row 7
$ws.Range("BF7").Value = 651651.008
$ws.Range("BG7").Value = 786995.008
$ws.Range("BH7").Value = 1067230.976
# Row 8
$ws.Range("BF8").Value = 0
$ws.Range("BG8").Value = 0
$ws.Range("BH8").Value = 0
# Row 9
$ws.Range("BF9").Value = 347657.984
$ws.Range("BG9").Value = 365340
$ws.Range("BH9").Value = 618806.976
# Row 10
$ws.Range("BF10").Value = 22129
$ws.Range("BG10").Value = 17694
$ws.Range("BH10").Value = 21857
# Row 11
$ws.Range("BF11").Value = 194542
$ws.Range("BG11").Value = 158548.992
$ws.Range("BH11").Value = 300097.984
# Row 12
$ws.Range("BF12").Value = 631700.992
$ws.Range("BG12").Value = 352604
$ws.Range("BH12").Value = 42123
# Row 13
$ws.Range("BF13").Value = 0
$ws.Range("BG13").Value = 0
$ws.Range("BH13").Value = 0
# Row 14
$ws.Range("BF14").Value = 0
$ws.Range("BG14").Value = 0
$ws.Range("BH14").Value = 0
# Row 15
$ws.Range("BF15").Value = 0
$ws.Range("BG15").Value = 0
$ws.Range("BH15").Value = 0
# Row 16
$ws.Range("BF16").Value = 0
$ws.Range("BG16").Value = 0
$ws.Range("BH16").Value = 0
# Row 17
$ws.Range("BF17").Value = 0
$ws.Range("BG17").Value = 0
$ws.Range("BH17").Value = 0
# Row 18
$ws.Range("BF18").Value = 0
$ws.Range("BG18").Value = 0
$ws.Range("BH18").Value = 0
# Row 19
$ws.Range("BF19").Value = 516582.016
$ws.Range("BG19").Value = 226895.008
$ws.Range("BH19").Value = 0
# Row 20
$ws.Range("BF20").Value = 0
$ws.Range("BG20").Value = 0
$ws.Range("BH20").Value = 0
# Row 21
$ws.Range("BF21").Value = 0
$ws.Range("BG21").Value = 0
$ws.Range("BH21").Value = 0
# Row 22
$ws.Range("BF22").Value = 0
$ws.Range("BG22").Value = 0
$ws.Range("BH22").Value = 0
# Row 23
$ws.Range("BF23").Value = 11469478.912
$ws.Range("BG23").Value = 12714584.064
$ws.Range("BH23").Value = 14926721.024
# Row 24
$ws.Range("BF24").Value = 10914878.464
$ws.Range("BG24").Value = 10921526.272
$ws.Range("BH24").Value = 11775377.408
# Row 25
$ws.Range("BF25").Value = 0
$ws.Range("BG25").Value = 0
$ws.Range("BH25").Value = 0
# Row 26
$ws.Range("BF26").Value = 28310931.456
$ws.Range("BG26").Value = 30779117.568
$ws.Range("BH26").Value = 37293928.448
# Row 27
$ws.Range("BF27").Value = 3463811.072
$ws.Range("BG27").Value = 3634616.064
$ws.Range("BH27").Value = 5183664.128
# Row 28
$ws.Range("BF28").Value = 266892.992
$ws.Range("BG28").Value = 179616.992
$ws.Range("BH28").Value = 196040
# Row 29
$ws.Range("BF29").Value = 834777.9840000001
$ws.Range("BG29").Value = 636779.008
$ws.Range("BH29").Value = 1123725.952
# Row 30
$ws.Range("BF30").Value = 593006.976
$ws.Range("BG30").Value = 410552.992
$ws.Range("BH30").Value = 698729.9840000001
# Row 31
$ws.Range("BF31").Value = 1476563.968
$ws.Range("BG31").Value = 1571096.96
$ws.Range("BH31").Value = 1595463.04
# Row 32
$ws.Range("BF32").Value = 0
$ws.Range("BG32").Value = 0
$ws.Range("BH32").Value = 0
# Row 33
$ws.Range("BF33").Value = 0
$ws.Range("BG33").Value = 0
$ws.Range("BH33").Value = 0
# Row 34
$ws.Range("BF34").Value = 292568.992
$ws.Range("BG34").Value = 836569.9840000001
$ws.Range("BH34").Value = 1569704.96
# Row 35
$ws.Range("BF35").Value = 0
$ws.Range("BG35").Value = 0
$ws.Range("BH35").Value = 0
# Row 36
$ws.Range("BF36").Value = 0
$ws.Range("BG36").Value = 0
$ws.Range("BH36").Value = 0
# Row 37
$ws.Range("BF37").Value = 10968841.216
$ws.Range("BG37").Value = 11747972.096
$ws.Range("BH37").Value = 14343484.416
# Row 38
$ws.Range("BF38").Value = 8077837.824
$ws.Range("BG38").Value = 9523656.704
$ws.Range("BH38").Value = 11662388.224
# Row 39
$ws.Range("BF39").Value = 0
$ws.Range("BG39").Value = 0
$ws.Range("BH39").Value = 0
# Row 40
$ws.Range("BF40").Value = 823473.024
$ws.Range("BG40").Value = 38421
$ws.Range("BH40").Value = 42567
# Row 41
$ws.Range("BF41").Value = 0
$ws.Range("BG41").Value = 0
$ws.Range("BH41").Value = 572934.976
# Row 42
$ws.Range("BF42").Value = 0
$ws.Range("BG42").Value = 0
$ws.Range("BH42").Value = 0
# Row 43
$ws.Range("BF43").Value = 2067529.984
$ws.Range("BG43").Value = 2185893.888
$ws.Range("BH43").Value = 2065593.984
# Row 44
$ws.Range("BF44").Value = 0
$ws.Range("BG44").Value = 0
$ws.Range("BH44").Value = 0
# Row 45
$ws.Range("BF45").Value = 0
$ws.Range("BG45").Value = 0
$ws.Range("BH45").Value = 0
# Row 46
$ws.Range("BF46").Value = 0
$ws.Range("BG46").Value = 0
$ws.Range("BH46").Value = 0
# Row 47
$ws.Range("BF47").Value = 13878280.192
$ws.Range("BG47").Value = 15396530.176
$ws.Range("BH47").Value = 17766782.976
# Row 48
$ws.Range("BF48").Value = 5352792.064
$ws.Range("BG48").Value = 7611386.88
$ws.Range("BH48").Value = 7611386.88
# Row 49
$ws.Range("BF49").Value = -171876.992
$ws.Range("BG49").Value = -242064
$ws.Range("BH49").Value = -326651.008
# Row 50
$ws.Range("BF50").Value = 0
$ws.Range("BG50").Value = 0
$ws.Range("BH50").Value = 0
# Row 51
$ws.Range("BF51").Value = 8801740.800000001
$ws.Range("BG51").Value = 6601740.8
$ws.Range("BH51").Value = 6601740.8
# Row 52
$ws.Range("BF52").Value = 0
$ws.Range("BG52").Value = 1045553.024
$ws.Range("BH52").Value = 2465562.112
# Row 53
$ws.Range("BF53").Value = 80284
$ws.Range("BG53").Value = 80526
$ws.Range("BH53").Value = -634193.9840000001
# Row 54
$ws.Range("BF54").Value = -184660
$ws.Range("BG54").Value = 299387.008
$ws.Range("BH54").Value = 2048936.96
# Row 55
$ws.Range("BF55").Value = 0
$ws.Range("BG55").Value = 0
$ws.Range("BH55").Value = 0
# Row 56
$ws.Range("BF56").Value = 0
$ws.Range("BG56").Value = 0
$ws.Range("BH56").Value = 0
# Row 59
$ws.Range("BF59").Value = 2733390.592
$ws.Range("BG59").Value = 3200186.112
$ws.Range("BH59").Value = 4595040.768
# Row 60
$ws.Range("BF60").Value = -793581.12
$ws.Range("BG60").Value = -1387639.04
$ws.Range("BH60").Value = -2190580.992
# Row 61
$ws.Range("BF61").Value = 1939808.896
$ws.Range("BG61").Value = 1812546.944
$ws.Range("BH61").Value = 2404460.032
# Row 62
$ws.Range("BF62").Value = 261808
$ws.Range("BG62").Value = 0
$ws.Range("BH62").Value = 0
# Row 63
$ws.Range("BF63").Value = -246839.024
$ws.Range("BG63").Value = -118814
$ws.Range("BH63").Value = -193616.992
# Row 64
$ws.Range("BF64").Value = 0
$ws.Range("BG64").Value = 0
$ws.Range("BH64").Value = 0
# Row 65
$ws.Range("BF65").Value = 64351
$ws.Range("BG65").Value = 30023
$ws.Range("BH65").Value = 309233.984
# Row 66
$ws.Range("BF66").Value = 203376.992
$ws.Range("BG66").Value = 0
$ws.Range("BH66").Value = 0
# Row 67
$ws.Range("BF67").Value = 0
$ws.Range("BG67").Value = 0
$ws.Range("BH67").Value = 0
# Row 68
$ws.Range("BF68").Value = -333748.992
$ws.Range("BG68").Value = -198679.008
$ws.Range("BH68").Value = -55870
# Row 69
$ws.Range("BF69").Value = 749409.9840000001
$ws.Range("BG69").Value = 402193.984
$ws.Range("BH69").Value = 2293392.896
# Row 70
$ws.Range("BF70").Value = -1083159.168
$ws.Range("BG70").Value = -600873.024
$ws.Range("BH70").Value = -2349263.104
# Row 74
$ws.Range("BF74").Value = 1888757.248
$ws.Range("BG74").Value = 1525076.992
$ws.Range("BH74").Value = 2464207.104
# Row 75
$ws.Range("BF75").Value = -262941.024
$ws.Range("BG75").Value = -180923.008
$ws.Range("BH75").Value = -244863.008
# Row 76
$ws.Range("BF76").Value = -77768.992
$ws.Range("BG76").Value = -298600.992
$ws.Range("BH76").Value = -799334.976
# Row 79
$ws.Range("BF79").Value = 0
$ws.Range("BG79").Value = 0
$ws.Range("BH79").Value = 0
# Row 80
$ws.Range("BF80").Value = 1548046.848
$ws.Range("BG80").Value = 1045553.024
$ws.Range("BH80").Value = 1420008.96
